$d = $word.ActiveDocument
$CR = [char]13

function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd($CR)
}

# --- 1. Extend the 04/11/2024 entry -----------------------------------
# Locate the paragraph "Started work on bandit problem proof of concept."
# and grow the trailing "." run into the fuller sentence, then append a
# brand-new run for the "Set up meeting..." sentence - all without
# disturbing the pre-existing "Started work on..." run.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Get-ParaText $d.Paragraphs.Item($i)) -eq "Started work on bandit problem proof of concept.") {
        $targetIndex = $i
    }
}

$pr = $d.Paragraphs.Item($targetIndex).Range
$periodRange = $d.Range($pr.End - 2, $pr.End - 1)
$periodRange.InsertAfter(" Started the report and got a good foundation. Created a rough UML for the program to ensure I planned it out properly. ")

$pr2 = $d.Paragraphs.Item($targetIndex).Range
$endPoint = $d.Range($pr2.End - 1, $pr2.End - 1)
$endPoint.InsertAfter("Set up meeting with supervisor for tomorrow. ")

# --- 2. Insert the 05/11/2024 entry between the two blank paragraphs --
# that follow the 04/11/2024 entry (the first of the two blank
# paragraphs immediately follows $targetIndex).
$blankIndex = $targetIndex + 1

$d.Paragraphs.Item($blankIndex).Range.InsertParagraphAfter()
$dateIndex = $blankIndex + 1
$d.Paragraphs.Item($dateIndex).Range.Text = "05/11/2024"

$d.Paragraphs.Item($dateIndex).Range.InsertParagraphAfter()
$meetingIndex = $dateIndex + 1
$d.Paragraphs.Item($meetingIndex).Range.Text = "Had meeting with supervisor"

# --- 3. Append the 06/11/2024 entry at the very end of the document ---
$newDatePara = $d.Paragraphs.Add()
$newDatePara.Range.Text = "06/11/2024"

$newBodyPara = $d.Paragraphs.Add()
$newBodyPara.Range.Text = "Continue work on UML and code for bandit problem proof of concept"
